$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add a thin paragraph border (space-only, no line) to the first paragraph.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.LeftIndent = 11.25

# Update the ID placeholder text.
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5345_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5345__ID**", 2)

# Drop the now-orphaned trailing space run that followed the placeholder.
$pEnd = $p1.Range.End
$spaceRange = $d.Range($pEnd - 2, $pEnd - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Text = ""
}
